# Update '想去人数' (attendance interest count) figures in column F
# across the 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types)
# worksheets, matching the refreshed data snapshot committed upstream.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
# Row 3: 杭州·COMIC TIME动漫游戏嘉年华 -> F3 3360 => 3364
$ws.Cells.Item(3, 6).Value = 3364
# Row 4: 杭州·异次结界动漫嘉年华 -> F4 1680 => 1679
$ws.Cells.Item(4, 6).Value = 1679
# Row 5: 杭州·ComicMe动漫嘉年华 -> F5 2422 => 2426
$ws.Cells.Item(5, 6).Value = 2426
# Row 7: 杭州·OZ·富坚义博only -> F7 338 => 339
$ws.Cells.Item(7, 6).Value = 339
# Row 8: 杭州·SST动漫嘉年华 -> F8 1373 => 1376
$ws.Cells.Item(8, 6).Value = 1376
# Row 9: 杭州·排球少年*蓝锁ONLY -> F9 1086 => 1089
$ws.Cells.Item(9, 6).Value = 1089
# Row 10: 杭州·春和景明代号鸢only -> F10 294 => 295
$ws.Cells.Item(10, 6).Value = 295
# Row 11: 杭州·百鬼夜行·咒术回战only -> F11 508 => 510
$ws.Cells.Item(11, 6).Value = 510
# Row 12: 杭州·第十届次元鹿角动漫游戏展（取消） -> F12 1166 => 1165
$ws.Cells.Item(12, 6).Value = 1165
# Row 16: 杭州·AD02动漫展 -> F16 8472 => 8501
$ws.Cells.Item(16, 6).Value = 8501
# Row 17: 杭州·AD02动漫展  青柳尊哉内场票 -> F17 370 => 371
$ws.Cells.Item(17, 6).Value = 371
# Row 18: 杭州·AD02动漫展--卡琳娜签售票 -> F18 2480 => 2479
$ws.Cells.Item(18, 6).Value = 2479
# Row 19: 杭州·现世繁华-代号鸢only -> F19 250 => 251
$ws.Cells.Item(19, 6).Value = 251
# Row 21: 杭州·AD02动漫展--亦之紫F、L句号内场票 -> F21 172 => 173
$ws.Cells.Item(21, 6).Value = 173
# Row 23: 杭州·Look Look动漫嘉年华 -> F23 580 => 582
$ws.Cells.Item(23, 6).Value = 582
# Row 27: 杭州·2024ESCC游戏电竞博览会暨新次元微光青春动漫交流会 -> F27 1994 => 2005
$ws.Cells.Item(27, 6).Value = 2005
# Row 28: 杭州·ELECTRIC COMIC动漫游戏展 -> F28 2034 => 2041
$ws.Cells.Item(28, 6).Value = 2041
# Row 30: 杭州·梦漫星河动漫展 -> F30 1731 => 1737
$ws.Cells.Item(30, 6).Value = 1737
# Row 32: 杭州·ESCC电竞博览会 倒霉死勒内场票 -> F32 1914 => 1915
$ws.Cells.Item(32, 6).Value = 1915
# Row 34: 杭州·第36届 中二病 原神x星穹only -> F34 25 => 26
$ws.Cells.Item(34, 6).Value = 26
# Row 37: 杭州·赛马娘only—晴空雏菊 -> F37 179 => 180
$ws.Cells.Item(37, 6).Value = 180
# Row 39: 杭州·SK怀旧展&偶像专场 -> F39 301 => 302
$ws.Cells.Item(39, 6).Value = 302
# Row 40: 杭州·【海潮的回响Echo of The Tide】 | 刀客塔们的大群融入派对·明日方舟SPECIAL ONLY -> F40 56 => 57
$ws.Cells.Item(40, 6).Value = 57
# Row 41: 杭州·白日梦次元动漫嘉年华 -> F41 232 => 235
$ws.Cells.Item(41, 6).Value = 235
# Row 43: 杭州·第四届华盟动漫次元嘉年华 -> F43 116 => 124
$ws.Cells.Item(43, 6).Value = 124
# Row 45: 杭州·次元幻想--二次元全女夜场 -> F45 253 => 255
$ws.Cells.Item(45, 6).Value = 255

$ws = $wb.Worksheets.Item("演出")
# Row 3: 杭州·《天空之城》久石让·宫崎骏动漫经典作品音乐会|浙江电影爱乐乐团 -> F3 14 => 15
$ws.Cells.Item(3, 6).Value = 15

$ws = $wb.Worksheets.Item("全部类型")
# Row 3: 杭州·COMIC TIME动漫游戏嘉年华 -> F3 3360 => 3364
$ws.Cells.Item(3, 6).Value = 3364
# Row 4: 杭州·异次结界动漫嘉年华 -> F4 1680 => 1679
$ws.Cells.Item(4, 6).Value = 1679
# Row 5: 杭州·ComicMe动漫嘉年华 -> F5 2422 => 2426
$ws.Cells.Item(5, 6).Value = 2426
# Row 7: 杭州·OZ·富坚义博only -> F7 338 => 339
$ws.Cells.Item(7, 6).Value = 339
# Row 8: 杭州·SST动漫嘉年华 -> F8 1373 => 1376
$ws.Cells.Item(8, 6).Value = 1376
# Row 10: 杭州·排球少年*蓝锁ONLY -> F10 1086 => 1089
$ws.Cells.Item(10, 6).Value = 1089
# Row 11: 杭州·春和景明代号鸢only -> F11 294 => 295
$ws.Cells.Item(11, 6).Value = 295
# Row 12: 杭州·百鬼夜行·咒术回战only -> F12 508 => 510
$ws.Cells.Item(12, 6).Value = 510
# Row 13: 杭州·第十届次元鹿角动漫游戏展（取消） -> F13 1166 => 1165
$ws.Cells.Item(13, 6).Value = 1165
# Row 16: 杭州·AD02动漫展 -> F16 8472 => 8501
$ws.Cells.Item(16, 6).Value = 8501
# Row 17: 杭州·AD02动漫展  青柳尊哉内场票 -> F17 370 => 371
$ws.Cells.Item(17, 6).Value = 371
# Row 18: 杭州·AD02动漫展--卡琳娜签售票 -> F18 2480 => 2479
$ws.Cells.Item(18, 6).Value = 2479
# Row 19: 杭州·《天空之城》久石让·宫崎骏动漫经典作品音乐会|浙江电影爱乐乐团 -> F19 14 => 15
$ws.Cells.Item(19, 6).Value = 15
# Row 20: 杭州·现世繁华-代号鸢only -> F20 250 => 251
$ws.Cells.Item(20, 6).Value = 251
# Row 22: 杭州·AD02动漫展--亦之紫F、L句号内场票 -> F22 172 => 173
$ws.Cells.Item(22, 6).Value = 173
# Row 24: 杭州·Look Look动漫嘉年华 -> F24 580 => 582
$ws.Cells.Item(24, 6).Value = 582
# Row 28: 杭州·2024ESCC游戏电竞博览会暨新次元微光青春动漫交流会 -> F28 1994 => 2005
$ws.Cells.Item(28, 6).Value = 2005
# Row 29: 杭州·ELECTRIC COMIC动漫游戏展 -> F29 2034 => 2041
$ws.Cells.Item(29, 6).Value = 2041
# Row 30: 杭州·梦漫星河动漫展 -> F30 1731 => 1737
$ws.Cells.Item(30, 6).Value = 1737
# Row 32: 杭州·ESCC电竞博览会 倒霉死勒内场票 -> F32 1914 => 1915
$ws.Cells.Item(32, 6).Value = 1915
# Row 34: 杭州·第36届 中二病 原神x星穹only -> F34 25 => 26
$ws.Cells.Item(34, 6).Value = 26
# Row 37: 杭州·赛马娘only—晴空雏菊 -> F37 179 => 180
$ws.Cells.Item(37, 6).Value = 180
# Row 39: 杭州·SK怀旧展&偶像专场 -> F39 301 => 302
$ws.Cells.Item(39, 6).Value = 302
# Row 40: 杭州·【海潮的回响Echo of The Tide】 | 刀客塔们的大群融入派对·明日方舟SPECIAL ONLY -> F40 56 => 57
$ws.Cells.Item(40, 6).Value = 57
# Row 41: 杭州·白日梦次元动漫嘉年华 -> F41 232 => 235
$ws.Cells.Item(41, 6).Value = 235
# Row 47: 杭州·第四届华盟动漫次元嘉年华 -> F47 116 => 124
$ws.Cells.Item(47, 6).Value = 124
# Row 49: 杭州·次元幻想--二次元全女夜场 -> F49 253 => 255
$ws.Cells.Item(49, 6).Value = 255
